$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (TC_02): "Account Id should be of 12 digits." -> "Account number should be of 8 digits"
# and the related error message changes from 12 digits to 8 digits.
$ws.Range("G6").Value = 'Error message:"Account Id should be of 8 digits"'
$ws.Range("D6").Value = "Account number should be of 8 digits"

# Row 18 (TC_14): same text change as row 6
$ws.Range("G18").Value = 'Error message:"Account Id should be of 8 digits"'
$ws.Range("D18").Value = "Account number should be of 8 digits"

# Row 19 (TC_15): same text change, and its G column (blank placeholder) becomes the error message
$ws.Range("G19").Value = 'Error message:"Account Id should be of 8 digits"'
$ws.Range("D19").Value = "Account number should be of 8 digits"

# Row 20 (TC_16): "Adhaar number should be of 12 digits" -> "Account number should be of 8 digits",
# and its G column (blank placeholder) becomes the error message
$ws.Range("G20").Value = 'Error message:"Account Id should be of 8 digits"'
$ws.Range("D20").Value = "Account number should be of 8 digits"

# Update the active selection on the sheet to D6
$ws.Range("D6").Select()
